$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.358.24'
$ws.Range('E2').Value = '  +1.96%  '
$ws.Range('D3').Value = '2.642.82'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '598.68'
$ws.Range('E5').Value = '  +1.53%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.14'
$ws.Range('E6').Value = '  +2.69%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.545'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '2.640.25'
$ws.Range('E9').Value = '  +1.52%  '
$ws.Range('E10').Value = '  +7.42%  '
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('E13').Value = '  +1.64%  '
$ws.Range('E14').Value = '  +2.68%  '
$ws.Range('E15').Value = '  +3.11%  '
$ws.Range('D16').Value = '3.128.41'
$ws.Range('E16').Value = '  +1.74%  '
$ws.Range('D17').Value = '68.307.70'
$ws.Range('D18').Value = '2.644.84'
$ws.Range('E18').Value = '  +1.89%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.38'
$ws.Range('E19').Value = '  +3.38%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '365.19'
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.38'
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.25'
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.84'
$ws.Range('E23').Value = '  +2.35%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.10'
$ws.Range('E24').Value = '  +2.62%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '73.26'
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.97'
$ws.Range('E27').Value = '  +0.55%  '
$ws.Range('D28').Value = '2.777.13'
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('E29').Value = '  +5.42%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '571.70'
$ws.Range('E31').Value = '  -1.72%  '
$ws.Range('E32').Value = '  +4.77%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.97'
$ws.Range('E33').Value = '  +4.25%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.85'
$ws.Range('E34').Value = '  +2.70%  '
$ws.Range('E35').Value = '  +2.95%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  +3.50%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '160.19'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '19.22'
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('E40').Value = '  +4.37%  '
$ws.Range('E41').Value = '  +0.82%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.37'
$ws.Range('E42').Value = '  +2.89%  '
$ws.Range('E43').Value = '  +3.74%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.62'
$ws.Range('E44').Value = '  +2.35%  '
$ws.Range('E45').Value = '  +12.25%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '40.53'
$ws.Range('E47').Value = '  -0.32%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '156.91'
$ws.Range('E48').Value = '  +2.69%  '
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('E50').Value = '  +1.68%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '21.82'
$ws.Range('E51').Value = '  +2.45%  '